# Update the Bento multi-filter query cells from a "grouped_recurrence_score"
# of "16-20" to "31-35" (per commit: "updated bento scripts as per
# availability of objects for new data set"), then restore the previously
# selected/active cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldScore = "16-20"
$newScore = "31-35"

# Columns B and C on rows 2-4 hold the Cypher query text for the
# CasesTab / SamplesTab / FilesTab rows (column C always repeats the basic
# stat-count query). Replace the recurrence-score filter in each one, in
# the same order Excel originally touched them so the shared-string table
# is rebuilt in the same order as the saved workbook.
foreach ($addr in @("B4", "C2", "C3", "C4", "B2", "B3")) {
    $cell = $ws.Range($addr)
    $text = $cell.Value2
    $cell.Value = $text.Replace($oldScore, $newScore)
}

# Move the active selection from D4 to B3, and scroll so row 3 is the
# top-left visible row (matches the saved view in the edited workbook).
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("B3").Select()
